$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62, shifting existing rows 62:126 down to 63:127.
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with its data.
$ws.Range("A62").Value = 11
$ws.Range("B62").Value = "Vega Monumental Concepción"
$ws.Range("C62").Value = "Bíobío"
$ws.Range("D62").Value = 44792
$ws.Range("E62").Value = 8
$ws.Range("F62").Value = 100112021
$ws.Range("G62").Value = "Ají"
$ws.Range("H62").Value = "Inferno"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 260
$ws.Range("K62").Value = 16000
$ws.Range("L62").Value = 18000
$ws.Range("M62").Value = 16923
$ws.Range("N62").Value = "$/caja 15 kilos"
$ws.Range("O62").Value = "Provincia de Huasco"
$ws.Range("P62").Value = 1128
$ws.Range("Q62").Value = 15
$ws.Range("R62").Value = "Hortaliza"
